$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 9).Value = "b"
$ws.Cells.Item(2, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(9, 9).Value = "sd"
$ws.Cells.Item(9, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(16, 9).Value = "sd"
$ws.Cells.Item(16, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(22, 9).Value = "aa"
$ws.Cells.Item(22, 10).Value = "Agree/Accept"
$ws.Cells.Item(25, 9).Value = "sd"
$ws.Cells.Item(25, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(37, 9).Value = "aa"
$ws.Cells.Item(37, 10).Value = "Agree/Accept"
$ws.Cells.Item(39, 9).Value = "sd"
$ws.Cells.Item(39, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(47, 9).Value = "aa"
$ws.Cells.Item(47, 10).Value = "Agree/Accept"
$ws.Cells.Item(49, 9).Value = "sv"
$ws.Cells.Item(49, 10).Value = "Statement-opinion"
$ws.Cells.Item(55, 9).Value = "sv"
$ws.Cells.Item(55, 10).Value = "Statement-opinion"
$ws.Cells.Item(59, 9).Value = "aa"
$ws.Cells.Item(59, 10).Value = "Agree/Accept"
$ws.Cells.Item(71, 9).Value = "aa"
$ws.Cells.Item(71, 10).Value = "Agree/Accept"
$ws.Cells.Item(90, 9).Value = "sv"
$ws.Cells.Item(90, 10).Value = "Statement-opinion"
$ws.Cells.Item(97, 9).Value = "%"
$ws.Cells.Item(97, 10).Value = "Uninterpretable"
$ws.Cells.Item(99, 9).Value = "aa"
$ws.Cells.Item(99, 10).Value = "Agree/Accept"
$ws.Cells.Item(104, 9).Value = "sv"
$ws.Cells.Item(104, 10).Value = "Statement-opinion"
$ws.Cells.Item(111, 9).Value = "%"
$ws.Cells.Item(111, 10).Value = "Uninterpretable"
$ws.Cells.Item(113, 9).Value = "sv"
$ws.Cells.Item(113, 10).Value = "Statement-opinion"
$ws.Cells.Item(115, 9).Value = "sd"
$ws.Cells.Item(115, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(116, 9).Value = "%"
$ws.Cells.Item(116, 10).Value = "Uninterpretable"
$ws.Cells.Item(125, 9).Value = "ba"
$ws.Cells.Item(125, 10).Value = "Appreciation"
$ws.Cells.Item(130, 9).Value = "sd"
$ws.Cells.Item(130, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(135, 9).Value = "b"
$ws.Cells.Item(135, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(136, 9).Value = "aa"
$ws.Cells.Item(136, 10).Value = "Agree/Accept"
$ws.Cells.Item(140, 9).Value = "ba"
$ws.Cells.Item(140, 10).Value = "Appreciation"
$ws.Cells.Item(145, 9).Value = "sd"
$ws.Cells.Item(145, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(160, 9).Value = "sd"
$ws.Cells.Item(160, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(161, 9).Value = "sd"
$ws.Cells.Item(161, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(162, 9).Value = "sd"
$ws.Cells.Item(162, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(170, 9).Value = "ba"
$ws.Cells.Item(170, 10).Value = "Appreciation"
$ws.Cells.Item(185, 9).Value = "ba"
$ws.Cells.Item(185, 10).Value = "Appreciation"
$ws.Cells.Item(193, 9).Value = "b"
$ws.Cells.Item(193, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(198, 9).Value = "sv"
$ws.Cells.Item(198, 10).Value = "Statement-opinion"
$ws.Cells.Item(199, 9).Value = "%"
$ws.Cells.Item(199, 10).Value = "Uninterpretable"
$ws.Cells.Item(233, 9).Value = "ba"
$ws.Cells.Item(233, 10).Value = "Appreciation"
$ws.Cells.Item(282, 9).Value = "%"
$ws.Cells.Item(282, 10).Value = "Uninterpretable"
$ws.Cells.Item(290, 9).Value = "%"
$ws.Cells.Item(290, 10).Value = "Uninterpretable"
$ws.Cells.Item(292, 9).Value = "ba"
$ws.Cells.Item(292, 10).Value = "Appreciation"
$ws.Cells.Item(296, 9).Value = "%"
$ws.Cells.Item(296, 10).Value = "Uninterpretable"
$ws.Cells.Item(306, 9).Value = "sv"
$ws.Cells.Item(306, 10).Value = "Statement-opinion"
$ws.Cells.Item(314, 9).Value = "ba"
$ws.Cells.Item(314, 10).Value = "Appreciation"
$ws.Cells.Item(315, 9).Value = "b"
$ws.Cells.Item(315, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(318, 9).Value = "sv"
$ws.Cells.Item(318, 10).Value = "Statement-opinion"
$ws.Cells.Item(322, 9).Value = "%"
$ws.Cells.Item(322, 10).Value = "Uninterpretable"
$ws.Cells.Item(324, 9).Value = "sv"
$ws.Cells.Item(324, 10).Value = "Statement-opinion"
$ws.Cells.Item(343, 9).Value = "sv"
$ws.Cells.Item(343, 10).Value = "Statement-opinion"
$ws.Cells.Item(351, 9).Value = "aa"
$ws.Cells.Item(351, 10).Value = "Agree/Accept"
